$wb = $excel.ActiveWorkbook

# XlLineStyle / XlBordersIndex constants used below:
#   xlContinuous = 1
#   xlEdgeTop = 8 ; xlEdgeBottom = 9 ; xlEdgeRight = 10

function Set-TopBottomBorderStyle($range) {
    # Reset to the plain default style first, then add only top+bottom thin
    # borders -> borderId 4 (matches the new cellXfs style index 2).
    $range.Style = "Normal"
    $range.Borders.Item(8).LineStyle = 1
    $range.Borders.Item(9).LineStyle = 1
}

function Set-RightTopBottomBorderStyle($range) {
    # Reset to the plain default style first, then add right+top+bottom thin
    # borders -> borderId 5 (matches the new cellXfs style index 3). Setting
    # the right edge first avoids materialising a transient "top only"
    # cellXfs entry that would otherwise linger unused.
    $range.Style = "Normal"
    $range.Borders.Item(10).LineStyle = 1
    $range.Borders.Item(8).LineStyle = 1
    $range.Borders.Item(9).LineStyle = 1
}

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item("quality_comparison")

Set-TopBottomBorderStyle $ws1.Range("C1")
Set-RightTopBottomBorderStyle $ws1.Range("D1")

$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item("computational_comparison")

Set-TopBottomBorderStyle $ws2.Range("C1")
Set-RightTopBottomBorderStyle $ws2.Range("D1")
Set-TopBottomBorderStyle $ws2.Range("F1")
Set-RightTopBottomBorderStyle $ws2.Range("G1")

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell at G5
$ws2.Range("G5").ClearContents()
